# A new weekly price record was added for "Betarraga" (Macroferia Regional
# de Talca). It belongs chronologically right before the existing row 157
# (old date 2021-05-25 / serial 44326), so insert a new row at position 157
# - this shifts the old rows 157..255 down to 158..256 - and fill the new
# row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 157, pushing rows 157..255 down to 158..256.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new record's data.
$ws.Cells.Item(157, 1).Value = 5
$ws.Cells.Item(157, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(157, 3).Value = "Maule"
$ws.Cells.Item(157, 4).Value = 44596
$ws.Cells.Item(157, 5).Value = 7
$ws.Cells.Item(157, 6).Value = 100114014
$ws.Cells.Item(157, 7).Value = "Betarraga"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 3000
$ws.Cells.Item(157, 11).Value = 700
$ws.Cells.Item(157, 12).Value = 700
$ws.Cells.Item(157, 13).Value = 700
$ws.Cells.Item(157, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 140
$ws.Cells.Item(157, 17).Value = 5
$ws.Cells.Item(157, 18).Value = "Hortaliza"
